$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.668.57'
$ws.Range("E2").Value = '  -2.49%  '
$ws.Range("D3").Value = '1.847.25'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4261'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3677'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.66'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07244'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9003'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.69'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.80%  '
$ws.Range("D13").Value = '1.864.77'
$ws.Range("E13").Value = '  -2.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.577'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.91%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.351'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06842'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.60%  '
$ws.Range("E17").Value = '  -0.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '77.62'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008820'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.42%  '
$ws.Range("D22").Value = '27.665.90'
$ws.Range("E22").Value = '  -2.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.961'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.09%  '
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("D25").Value = '2.089.08'
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.046'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.90'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.16'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.250'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.828'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '111.07'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08886'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7704'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.565'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.907'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.084'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.001'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.095'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05373'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.977'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01924'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5061'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.767'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1639'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.264'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06638'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.19%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4721'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.38%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.14'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.001'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.638'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.66%  '
